$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Menu category headers: uppercase wording + a refreshed (explicit black) font colour ---
# Row 10: Steak -> STEAK
$ws.Range("A10").Value = "STEAK"
# Row 12: Chicken -> CHICKEN
$ws.Range("A12").Value = "CHICKEN"
# Row 14: Pasta -> PASTA
$ws.Range("A14").Value = "PASTA"
# Row 16: was mistakenly "Wine" -> corrected to PIZZA
$ws.Range("A16").Value = "PIZZA"
# Row 18: Wine -> WINE
$ws.Range("A18").Value = "WINE"
# Row 20: Spirits -> SPIRITS
$ws.Range("A20").Value = "SPIRITS"

# Apply the refreshed font (Arial 16, explicit black) to the category cells and
# to the blank spacer cells directly beneath each of them.
$categoryRows = @(10,11,12,13,14,15,16,17,18,19,20)
foreach ($row in $categoryRows) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 16
    $cell.Font.Color = 0
}

# --- Item detail corrections ---
# Row 20 item description: Top Shelf -> Well
$ws.Range("C20").Value = "Well"

# --- Sales price corrections ---
$ws.Range("J18").Value = 7
$ws.Range("J20").Value = 6

# --- Last active selection when the author saved the file ---
$ws.Range("J21").Select()
